# Update the "想去人数" (interested-count) and "最低票价" (min price) figures
# across the four sheets, matching the upstream data refresh captured at
# commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 862
$ws.Range("F3").Value  = 13798
$ws.Range("F4").Value  = 13591
$ws.Range("F5").Value  = 1051
$ws.Range("F8").Value  = 597
$ws.Range("F12").Value = 764
$ws.Range("F13").Value = 2147
$ws.Range("F14").Value = 103
$ws.Range("F16").Value = 75
$ws.Range("F19").Value = 529
$ws.Range("F20").Value = 434
$ws.Range("F21").Value = 405
$ws.Range("F22").Value = 324
$ws.Range("F23").Value = 264
$ws.Range("F24").Value = 837
$ws.Range("F25").Value = 89
$ws.Range("F26").Value = 3

# --- Sheet "演出" -----------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 72
$ws.Range("G6").Value = 280
$ws.Range("F7").Value = 1513

# --- Sheet "本地生活" --------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 110

# --- Sheet "全部类型" --------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 862
$ws.Range("F4").Value  = 13798
$ws.Range("F5").Value  = 13591
$ws.Range("F6").Value  = 1051
$ws.Range("F9").Value  = 597
$ws.Range("F13").Value = 764
$ws.Range("F16").Value = 2147
$ws.Range("F17").Value = 103
$ws.Range("F19").Value = 75
$ws.Range("F23").Value = 72
$ws.Range("F24").Value = 110
$ws.Range("F25").Value = 110
$ws.Range("F26").Value = 529
$ws.Range("F27").Value = 434
$ws.Range("F28").Value = 405
$ws.Range("F29").Value = 324
$ws.Range("F30").Value = 264
$ws.Range("F31").Value = 837
$ws.Range("G32").Value = 280
$ws.Range("F33").Value = 1513
$ws.Range("F37").Value = 89
$ws.Range("F40").Value = 3
